$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.990.09"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.257.13"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.18"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.51"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.83%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0784"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.64%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.79"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.607.73"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.254.37"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.785"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.845.01"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0898"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.46%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.42"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.83"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.12%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.57"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.05%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.13"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.49"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.87"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.82%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.20"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.14"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0731"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.89%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.43"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.01%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.10"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.29"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.942.32"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.82"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.04%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.85"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.02%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.93"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.479.54"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "92.06"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.42"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.87%  "
